# The underlying source data re-ordered two species records (the one that
# was previously in row 2 is now in row 4, and vice versa). Apply this by
# swapping the values of the columns that differ between the two rows:
# Id, Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn,
# Auktor and Publik kommentar.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","E","F","G","H","AC")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr4 = "$col" + "4"

    $v2 = $ws.Range($addr2).Value()
    $v4 = $ws.Range($addr4).Value()

    $ws.Range($addr2).Value = $v4
    $ws.Range($addr4).Value = $v2
}

# Column L ("Kön") only has an (empty) cell on one of the two rows; move
# that empty placeholder cell from row 2 to row 4 along with the rest of
# the record, preserving its presence/absence rather than just its value.
$ws.Range("L2").Copy($ws.Range("L4"))
$ws.Range("L2").ClearContents()
